$wb = $excel.ActiveWorkbook

# --- department: page setup touched while reviewing tables ---
$wsDept = $wb.Worksheets.Item("department")
$wsDept.PageSetup.PaperSize = 9
$wsDept.PageSetup.Orientation = 1

# --- user_department: add primary key column (id / 主键) as first field ---
$wsUserDept = $wb.Worksheets.Item("user_department")
$wsUserDept.Rows.Item(2).Insert()
$wsUserDept.Range("A2").Value = "id"
$wsUserDept.Range("B2").Value = "主键"
$wsUserDept.PageSetup.PaperSize = 9
$wsUserDept.PageSetup.Orientation = 1

# --- user_project: add primary key column (id / 主键) as first field ---
$wsUserProj = $wb.Worksheets.Item("user_project")
$wsUserProj.Rows.Item(2).Insert()
$wsUserProj.Range("A2").Value = "id"
$wsUserProj.Range("B2").Value = "主键"

# restore the selection on each touched sheet to the row below its data
[void]$wsUserDept.Range("A7").Select()
[void]$wsUserProj.Range("A8").Select()

# leave focus on the "task" sheet (matches the workbook's saved active tab)
[void]$wb.Worksheets.Item("task").Select()
